$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.819.69'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '2.043.25'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''227.52'
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("D7").Value = '''59.73'
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '''0.377'
$ws.Range("E9").Value = '  -2.57%  '
$ws.Range("E10").Value = '  +2.62%  '
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("D12").Value = '2.344.89'
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("D13").Value = '''14.47'
$ws.Range("E13").Value = '  -1.20%  '
$ws.Range("D14").Value = '''21.05'
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").Value = '''5.49'
$ws.Range("E15").Value = '  +5.06%  '
$ws.Range("E16").Value = '  +1.62%  '
$ws.Range("D17").Value = '2.047.45'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").Value = '37.764.77'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = '''69.54'
$ws.Range("E19").Value = '  -0.43%  '
$ws.Range("E20").Value = '  -2.80%  '
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '''223.84'
$ws.Range("E22").Value = '  -0.76%  '
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("E25").Value = '  +3.59%  '
$ws.Range("D26").Value = '''169.59'
$ws.Range("E26").Value = '  +2.73%  '
$ws.Range("D27").Value = '''9.37'
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("E28").Value = '  -0.56%  '
$ws.Range("D29").Value = '''18.81'
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("E31").Value = '  -0.71%  '
$ws.Range("D32").Value = '''2.24'
$ws.Range("E32").Value = '  +9.24%  '
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("D34").Value = '''0.0603'
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").Value = '''6.52'
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("E37").Value = '  +3.94%  '
$ws.Range("E38").Value = '  +6.18%  '
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("E40").Value = '  +6.23%  '
$ws.Range("D41").Value = '1.526.89'
$ws.Range("E41").Value = '  -0.87%  '
$ws.Range("D42").Value = '''97.61'
$ws.Range("E42").Value = '  +0.70%  '
$ws.Range("E43").Value = '  -0.53%  '
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("D45").Value = '''0.0907'
$ws.Range("E45").Value = '  -1.53%  '
$ws.Range("D46").Value = '''4.17'
$ws.Range("E46").Value = '  +6.47%  '
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("D51").Value = '2.234.10'
$ws.Range("E51").Value = '  +0.43%  '
